$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 'culture_collection' column (AF) entirely. This shifts every
# subsequent column (and its header cell, cell comment, etc.) one position
# to the left, matching the diff where each comment's text moved up into
# the preceding cell and the final CD15 comment disappeared.
$ws.Columns("AF").Delete()
